$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the FilesTab Cypher query (cell B4): reorder the RETURN clause so that
# `Format` is returned right after `File Name` (before `File Type`), matching
# the corrected query used for case-file comparisons (icdc fix).
$newFilesQuery = "MATCH (f:file)-->(parent)`n" +
"WITH DISTINCT f, parent`n" +
"MATCH (diag:diagnosis)-->(c)`n" +
"OPTIONAL MATCH (f)-[*]->(samp:sample)-->(c)-->(s:study)`n" +
"MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
"WHERE s.clinical_study_designation IN ['UBC02']`n" +
"OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)`n" +
"WITH`n" +
"        f, parent, c, demo, diag, s, samp,`n" +
"        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n" +
"        toInteger(floor(log(f.file_size)/log(1024))) as i,`n" +
"        2 as precision`n" +
"WITH`n" +
"        f, parent, c, demo, diag, s, samp,`n" +
"        f.file_size /(1024^i) AS value, `n" +
"        10^precision AS factor,`n" +
"        units[i] as unit`n" +
"WITH    `n" +
"        f, parent, c, demo, diag, s, samp, unit,`n" +
"        round(factor * value)/factor AS size`n" +
"RETURN `n" +
"       coalesce(f.file_name, '') AS ``File Name``,`n" +
"       coalesce(f.file_format, '') AS ``Format``,`n" +
"        coalesce(f.file_type, '') AS ``File Type``,`n" +
"       CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n" +
"        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
"        coalesce(f.file_description, '') AS ``Description``,`n" +
"        coalesce(samp.sample_id, '') AS ``Sample ID``,`n" +
"        coalesce(c.case_id, '') AS ``Case ID``,`n" +
"        coalesce(demo.breed,'') AS Breed ,`n" +
"        coalesce(diag.disease_term,'') AS Diagnosis"

# Remember the current (fixed) row height so that assigning the longer
# query text below doesn't trigger Excel's row auto-fit behavior.
$row4Height = $ws.Rows.Item(4).RowHeight

$ws.Range("B4").Value = $newFilesQuery

# Restore the original row height (the row stays a fixed custom height).
$ws.Rows.Item(4).RowHeight = $row4Height

# Move the active selection from B2 to D4, as recorded in the saved view state.
$ws.Range("D4").Select()

$wb.Save()
